$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 25 label: Developer Integrity Proof ---
$ws.Range("A25").Value = "Developer Integrity Proof"

# --- Update Developer Integrity row (row 24) text content ---
$ws.Range("D24").Value = "Claims A New Library Similar To TensorFlow Will Be Developed, But There Are No Development For The Last 2 Years."

# --- Fill in row 25 hyperlink display text, in the same order the links were originally authored ---
$ws.Range("D25").Value = "https://devforum.roblox.com/t/neural-network-library-20/869557/126?u=myoriginsworkshop"
$ws.Range("E25").Value = "https://devforum.roblox.com/t/openml-machine-learning/3008664/4?u=myoriginsworkshop"
$ws.Range("F25").Value = "https://devforum.roblox.com/t/xentorch-neural-network-constructor/1201111/24?u=myoriginsworkshop"

$ws.Range("F24").Value = "Busy With Academics, So No New Features."

$ws.Range("G25").Value = "https://devforum.roblox.com/t/easyml-an-easy-way-to-use-machine-learning-in-your-roblox-games/3110013?u=myoriginsworkshop"

# --- Turn the row-25 URLs into real hyperlinks (this also applies the built-in Hyperlink style) ---
$ws.Hyperlinks.Add($ws.Range("D25"), "https://devforum.roblox.com/t/neural-network-library-20/869557/126?u=myoriginsworkshop") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E25"), "https://devforum.roblox.com/t/openml-machine-learning/3008664/4?u=myoriginsworkshop") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://devforum.roblox.com/t/xentorch-neural-network-constructor/1201111/24?u=myoriginsworkshop") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G25"), "https://devforum.roblox.com/t/easyml-an-easy-way-to-use-machine-learning-in-your-roblox-games/3110013?u=myoriginsworkshop") | Out-Null

$ws.Range("A25:G25").RowHeight = 60

$ws.Range("G24").Select()

Write-Host "Done"
